$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 84754
$ws.Range("I9").Value = 91650.164
$ws.Range("K9").Value = 91650.164
$ws.Range("M9").Value = -91481.164

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 491.45456
$ws.Range("I12").Value = 125.875
$ws.Range("K12").Value = 125.875
$ws.Range("M12").Value = 44.125

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 804.7619
$ws.Range("J19").Value = 846.8
$ws.Range("L19").Value = 846.8
$ws.Range("N19").Value = -1196.8

# ALC row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 888.8889
$ws.Range("I29").Value = 375
$ws.Range("K29").Value = 1125
$ws.Range("M29").Value = -844

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 392.42856
$ws.Range("I41").Value = 469.4
$ws.Range("J41").Value = 200
$ws.Range("K41").Value = 469.4
$ws.Range("L41").Value = 200
$ws.Range("M41").Value = -29.39999999999998
$ws.Range("N41").Value = -1080

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5553.6665
$ws.Range("I62").Value = 4764.6
$ws.Range("K62").Value = 4764.6
$ws.Range("M62").Value = -4140.6

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5553.6665
$ws.Range("I65").Value = 4764.6
$ws.Range("K65").Value = 23823
$ws.Range("M65").Value = -20703

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 428.83334
$ws.Range("I80").Value = 272.5
$ws.Range("J80").Value = 624.25
$ws.Range("K80").Value = 817.5
$ws.Range("L80").Value = 1872.75
$ws.Range("M80").Value = 180.5
$ws.Range("N80").Value = -3868.75

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 428.83334
$ws.Range("I83").Value = 272.5
$ws.Range("J83").Value = 624.25
$ws.Range("K83").Value = 2452.5
$ws.Range("L83").Value = 5618.25
$ws.Range("M83").Value = 2539.5
$ws.Range("N83").Value = -15602.25

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1381.0869
$ws.Range("I137").Value = 1160.2632
$ws.Range("J137").Value = 2430
$ws.Range("K137").Value = 3480.7896
$ws.Range("L137").Value = 7290
$ws.Range("M137").Value = -930.7896000000001
$ws.Range("N137").Value = -12390

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1533.1578
$ws.Range("I2").Value = 1141
$ws.Range("J2").Value = 2205.4285
$ws.Range("K2").Value = 1141
$ws.Range("L2").Value = 2205.4285
$ws.Range("M2").Value = -1028
$ws.Range("N2").Value = -2431.4285

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4961.7373
$ws.Range("I32").Value = 3726.0144
$ws.Range("K32").Value = 3726.0144
$ws.Range("M32").Value = -3439.0144

# ARM row 39
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 13958
$ws.Range("I39").Value = 2916
$ws.Range("K39").Value = 2916
$ws.Range("M39").Value = -2396

# ARM row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 19933.666
$ws.Range("I41").Value = 9900
$ws.Range("K41").Value = 9900
$ws.Range("M41").Value = -9486

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 843.4545000000001
$ws.Range("I110").Value = 843.4545000000001
$ws.Range("K110").Value = 843.4545000000001
$ws.Range("M110").Value = 1201.5455

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1533.1578
$ws.Range("I116").Value = 1141
$ws.Range("J116").Value = 2205.4285
$ws.Range("K116").Value = 1141
$ws.Range("L116").Value = 2205.4285
$ws.Range("M116").Value = 1153
$ws.Range("N116").Value = -6793.4285

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1533.1578
$ws.Range("I3").Value = 1141
$ws.Range("J3").Value = 2205.4285
$ws.Range("K3").Value = 1141
$ws.Range("L3").Value = 2205.4285
$ws.Range("M3").Value = -1027
$ws.Range("N3").Value = -2433.4285

# CRP row 35
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 2278.5715
$ws.Range("I35").Value = 2490
$ws.Range("J35").Value = 1750
$ws.Range("K35").Value = 2490
$ws.Range("L35").Value = 1750
$ws.Range("M35").Value = -2196
$ws.Range("N35").Value = -2338

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 78.25
$ws.Range("J12").Value = 56.375
$ws.Range("L12").Value = 169.125
$ws.Range("N12").Value = -515.125

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 6089.75
$ws.Range("I70").Value = 1513
$ws.Range("J70").Value = 6330.6313
$ws.Range("K70").Value = 4539
$ws.Range("L70").Value = 18991.8939
$ws.Range("M70").Value = -4224
$ws.Range("N70").Value = -19621.8939

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 6089.75
$ws.Range("I73").Value = 1513
$ws.Range("J73").Value = 6330.6313
$ws.Range("K73").Value = 4539
$ws.Range("L73").Value = 18991.8939
$ws.Range("M73").Value = -3447
$ws.Range("N73").Value = -21175.8939

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 5557736
$ws.Range("I129").Value = 369.27274
$ws.Range("J129").Value = 8775159
$ws.Range("K129").Value = 1107.81822
$ws.Range("L129").Value = 26325477
$ws.Range("M129").Value = 3892.18178
$ws.Range("N129").Value = -26335477

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 17302.646
$ws.Range("I102").Value = 17302.646
$ws.Range("K102").Value = 17302.646
$ws.Range("M102").Value = -15680.646

# LTW row 13
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 8752
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3671
$ws.Range("I32").Value = 3671
$ws.Range("K32").Value = 3671
$ws.Range("M32").Value = -3354

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1195.0769
$ws.Range("I46").Value = 871
$ws.Range("K46").Value = 871
$ws.Range("M46").Value = -683

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1878.2963
$ws.Range("I61").Value = 1854.3846
$ws.Range("K61").Value = 1854.3846
$ws.Range("M61").Value = -1652.3846

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1878.2963
$ws.Range("I113").Value = 1854.3846
$ws.Range("K113").Value = 1854.3846
$ws.Range("M113").Value = 315.6153999999999

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3685.5312
$ws.Range("J122").Value = 3694.0557
$ws.Range("L122").Value = 11082.1671
$ws.Range("N122").Value = -15982.1671

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3082.48
$ws.Range("I136").Value = 2691.6296
$ws.Range("K136").Value = 8074.888800000001
$ws.Range("M136").Value = -5524.888800000001

# WVR row 33
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 17333.334
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2750

# WVR row 36
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H36").Value = 17333.334
$ws.Range("I36").Value = 3000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2750
